# AKK TAXI business card - add yellow Mercedes taxi illustration
# Applies the resize/reposition of 4 existing shapes and appends the
# 13 new shapes (id 8-20) that make up the stylized cab graphic.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1. Resize / reposition the existing shapes that shifted to make room
#    for the new illustration.
# ---------------------------------------------------------------------

# Shape id=3 "Rectangle 2" (rotated teal accent rectangle)
$shRect2 = $s.Shapes.Item(2)
$shRect2.Left   = 147.59992225984251
$shRect2.Top    = -18.0
$shRect2.Width  = 108.0
$shRect2.Height = 93.6000023

# Shape id=5 "TextBox 4" (Razvan Taxi / Servicii de transport 24/7) - only size changes
$shTb4 = $s.Shapes.Item(4)
$shTb4.Width  = 165.6
$shTb4.Height = 50.4

# Shape id=6 "TextBox 5" (Tel: 0720 064 963)
$shTb5 = $s.Shapes.Item(5)
$shTb5.Left   = 25.2
$shTb5.Top    = 108.0
$shTb5.Width  = 158.4000016
$shTb5.Height = 36.0

# Shape id=7 "TextBox 6" (Rapid / Sigur / Confortabil)
$shTb6 = $s.Shapes.Item(6)
$shTb6.Left   = 183.6
$shTb6.Top    = 28.8000002
$shTb6.Width  = 68.4
$shTb6.Height = 72.0

# ---------------------------------------------------------------------
# 2. Add the new shapes that make up the yellow Mercedes taxi graphic.
# ---------------------------------------------------------------------

# id=8 "TextBox 7" - "Mercedes Fleet" caption
$s8 = $s.Shapes.Item(3).Duplicate()
$s8.Name = "TextBox 7"
$s8.Left   = 183.6
$s8.Top    = 111.6000023
$s8.Width  = 68.4
$s8.Height = 25.2
$s8.TextFrame.TextRange.Text = "Mercedes Fleet"
$s8.TextFrame.TextRange.Font.Size = 10
$s8.TextFrame.TextRange.Font.Bold = $true
$s8.TextFrame.TextRange.Font.Color.RGB = 0x241A08

# id=9 "Rounded Rectangle 8" - yellow cab body (large)
$s9 = $s.Shapes.Item(1).Duplicate()
$s9.Name = "Rounded Rectangle 8"
$s9.AutoShapeType = 5
$s9.Left   = 18.0
$s9.Top    = 90.0
$s9.Width  = 115.2000008
$s9.Height = 32.4
$s9.Adjustments.Item(1) = 0.2
$s9.Fill.ForeColor.RGB = 0x20D0FC
$s9.Line.Visible = $false

# id=10 "Rounded Rectangle 9" - yellow cab roof/cabin
$s10 = $s.Shapes.Item(1).Duplicate()
$s10.Name = "Rounded Rectangle 9"
$s10.AutoShapeType = 5
$s10.Left   = 39.6000004
$s10.Top    = 75.6000023
$s10.Width  = 57.6000004
$s10.Height = 21.6
$s10.Adjustments.Item(1) = 0.4
$s10.Fill.ForeColor.RGB = 0x20D0FC
$s10.Line.Visible = $false

# id=11 "Rectangle 10" - dark windshield band
$s11 = $s.Shapes.Item(1).Duplicate()
$s11.Name = "Rectangle 10"
$s11.AutoShapeType = 1
$s11.Left   = 43.2
$s11.Top    = 79.2000008
$s11.Width  = 50.4
$s11.Height = 15.84
$s11.Fill.ForeColor.RGB = 0x5A3C1E
$s11.Line.Visible = $false

# id=12 "Oval 11" - left wheel (outer, dark)
$s12 = $s.Shapes.Item(1).Duplicate()
$s12.Name = "Oval 11"
$s12.AutoShapeType = 9
$s12.Left   = 32.4
$s12.Top    = 113.76
$s12.Width  = 25.2
$s12.Height = 25.2
$s12.Fill.ForeColor.RGB = 0x16110C
$s12.Line.Visible = $false

# id=13 "Oval 12" - right wheel (outer, dark)
$s13 = $s.Shapes.Item(1).Duplicate()
$s13.Name = "Oval 12"
$s13.AutoShapeType = 9
$s13.Left   = 90.0
$s13.Top    = 113.76
$s13.Width  = 25.2
$s13.Height = 25.2
$s13.Fill.ForeColor.RGB = 0x16110C
$s13.Line.Visible = $false

# id=14 "Oval 13" - left wheel hubcap
$s14 = $s.Shapes.Item(1).Duplicate()
$s14.Name = "Oval 13"
$s14.AutoShapeType = 9
$s14.Left   = 38.160001799999996
$s14.Top    = 119.5200005
$s14.Width  = 13.68
$s14.Height = 13.68
$s14.Fill.ForeColor.RGB = 0xBEBEBE
$s14.Line.Visible = $false

# id=15 "Oval 14" - right wheel hubcap
$s15 = $s.Shapes.Item(1).Duplicate()
$s15.Name = "Oval 14"
$s15.AutoShapeType = 9
$s15.Left   = 95.76
$s15.Top    = 119.5200005
$s15.Width  = 13.68
$s15.Height = 13.68
$s15.Fill.ForeColor.RGB = 0xBEBEBE
$s15.Line.Visible = $false

# id=16 "Rectangle 15" - taxi sign post
$s16 = $s.Shapes.Item(1).Duplicate()
$s16.Name = "Rectangle 15"
$s16.AutoShapeType = 1
$s16.Left   = 126.0
$s16.Top    = 95.04
$s16.Width  = 8.64
$s16.Height = 21.6
$s16.Fill.ForeColor.RGB = 0xF0F0F0
$s16.Line.Visible = $false

# id=17 "Oval 16" - taxi sign lamp
$s17 = $s.Shapes.Item(1).Duplicate()
$s17.Name = "Oval 16"
$s17.AutoShapeType = 9
$s17.Left   = 128.16
$s17.Top    = 101.5200005
$s17.Width  = 5.76
$s17.Height = 9.3600002
$s17.Fill.ForeColor.RGB = 0x96EBFF
$s17.Line.Visible = $false

# id=18 "Rectangle 17" - headlight / door highlight
$s18 = $s.Shapes.Item(1).Duplicate()
$s18.Name = "Rectangle 17"
$s18.AutoShapeType = 1
$s18.Left   = 66.2400017
$s18.Top    = 70.5600014
$s18.Width  = 28.8000002
$s18.Height = 8.64
$s18.Fill.ForeColor.RGB = 0xFFFFFF
$s18.Line.Visible = $false

# id=19 "TextBox 18" - "TAXI" caption on the roof sign
$s19 = $s.Shapes.Item(3).Duplicate()
$s19.Name = "TextBox 18"
$s19.Left   = 66.96000289999999
$s19.Top    = 69.12
$s19.Width  = 27.36
$s19.Height = 11.52
$s19.TextFrame.TextRange.Text = "TAXI"
$s19.TextFrame.TextRange.Font.Size = 9
$s19.TextFrame.TextRange.Font.Bold = $true
$s19.TextFrame.TextRange.Font.Color.RGB = 0x241A08

# id=20 "5-Point Star 19" - decorative star on the sign
$s20 = $s.Shapes.Item(1).Duplicate()
$s20.Name = "5-Point Star 19"
$s20.AutoShapeType = 92
$s20.Left   = 127.44
$s20.Top    = 98.6400033
$s20.Width  = 6.48
$s20.Height = 6.48
$s20.Fill.ForeColor.RGB = 0x241A08
$s20.Line.Visible = $false
